$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1): update column F (想去人数)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 771
$ws1.Range("F4").Value = 507
$ws1.Range("F5").Value = 255
$ws1.Range("F6").Value = 460
$ws1.Range("F7").Value = 1094
$ws1.Range("F8").Value = 314
$ws1.Range("F11").Value = 96
$ws1.Range("F12").Value = 1082
$ws1.Range("F15").Value = 723
$ws1.Range("F16").Value = 781
$ws1.Range("F18").Value = 24
$ws1.Range("F19").Value = 50
$ws1.Range("F20").Value = 615
$ws1.Range("F21").Value = 105
$ws1.Range("F22").Value = 1687
$ws1.Range("F23").Value = 1899
$ws1.Range("F24").Value = 480
$ws1.Range("F26").Value = 1724
$ws1.Range("F27").Value = 255
$ws1.Range("F28").Value = 2532
$ws1.Range("F29").Value = 447
$ws1.Range("F30").Value = 33
$ws1.Range("F31").Value = 649
$ws1.Range("F33").Value = 84
$ws1.Range("F35").Value = 871
$ws1.Range("F36").Value = 1581
$ws1.Range("F37").Value = 263
$ws1.Range("F39").Value = 509
$ws1.Range("F40").Value = 106
$ws1.Range("F41").Value = 96

# Sheet "演出" (sheet2): update column F (想去人数) and G (最低票价)
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("G2").Value = 108
$ws2.Range("F3").Value = 8
$ws2.Range("F4").Value = 119
$ws2.Range("F12").Value = 60

# Sheet "全部类型" (sheet4): update column F (想去人数) and G (最低票价)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("G3").Value = 108
$ws4.Range("F4").Value = 8
$ws4.Range("F5").Value = 771
$ws4.Range("F6").Value = 507
$ws4.Range("F7").Value = 255
$ws4.Range("F8").Value = 460
$ws4.Range("F9").Value = 1094
$ws4.Range("F10").Value = 314
$ws4.Range("F13").Value = 96
$ws4.Range("F14").Value = 1082
$ws4.Range("F16").Value = 723
$ws4.Range("F17").Value = 781
$ws4.Range("F19").Value = 119
$ws4.Range("F20").Value = 119
$ws4.Range("F22").Value = 24
$ws4.Range("F24").Value = 50
$ws4.Range("F25").Value = 615
$ws4.Range("F26").Value = 105
$ws4.Range("F27").Value = 1687
$ws4.Range("F28").Value = 1899
$ws4.Range("F29").Value = 480
$ws4.Range("F32").Value = 2532
$ws4.Range("F33").Value = 447
$ws4.Range("F36").Value = 33
$ws4.Range("F37").Value = 60
$ws4.Range("F38").Value = 649
$ws4.Range("F40").Value = 84
$ws4.Range("F42").Value = 871
$ws4.Range("F43").Value = 1581
$ws4.Range("F45").Value = 263
$ws4.Range("F46").Value = 509
$ws4.Range("F47").Value = 107
$ws4.Range("F48").Value = 96
